$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild rows 2-69 (A:D) reflecting the regrouped permeation residue combinations
# (individual permeation events re-sorted/re-grouped by frame, prior to adding distance check)
$data = New-Object 'object[,]' 68,4
$data[0,0] = "131, 456, 781, 781, 1074"
$data[0,1] = 1
$data[0,2] = "781"
$data[0,3] = "1661"
$data[1,0] = "131, 456, 781, SF"
$data[1,1] = 1
$data[1,2] = "781"
$data[1,3] = "235"
$data[2,0] = "99, 749, 1074, 1106"
$data[2,1] = 1
$data[2,2] = "1106"
$data[2,3] = "571"
$data[3,0] = "131, 131, 456, 1106"
$data[3,1] = 2
$data[3,2] = "131, 456"
$data[3,3] = "382, 3600"
$data[4,0] = "99, 131, 781"
$data[4,1] = 1
$data[4,2] = "131"
$data[4,3] = "458"
$data[5,0] = "131, 749, 781, 1074"
$data[5,1] = 9
$data[5,2] = "781, 781, 781, 781, 781, 781, 781, 781, 781"
$data[5,3] = "1091, 849, 1295, 2197, 2374, 2385, 2553, 2697, 2734"
$data[6,0] = "99, 456, 749, 781, 1074"
$data[6,1] = 1
$data[6,2] = "781"
$data[6,3] = "805"
$data[7,0] = "456, 749, 781"
$data[7,1] = 1
$data[7,2] = "781"
$data[7,3] = "1120"
$data[8,0] = "131, 749, 781, SF"
$data[8,1] = 1
$data[8,2] = "781"
$data[8,3] = "1716"
$data[9,0] = "99, 131, 749, 781, 1074"
$data[9,1] = 1
$data[9,2] = "781"
$data[9,3] = "1000"
$data[10,0] = "131, 456, 749, 781, 1074"
$data[10,1] = 3
$data[10,2] = "781, 781, 781"
$data[10,3] = "1416, 1331, 6643"
$data[11,0] = "99, 781, 1074, SF"
$data[11,1] = 1
$data[11,2] = "781"
$data[11,3] = "1489"
$data[12,0] = "99, 131, 749, 1074"
$data[12,1] = 2
$data[12,2] = "131, 131"
$data[12,3] = "1585, 2537"
$data[13,0] = "99, 781, 781, 1074"
$data[13,1] = 1
$data[13,2] = "781"
$data[13,3] = "1899"
$data[14,0] = "99, 456, 749, 1074, 1106"
$data[14,1] = 1
$data[14,2] = "1106"
$data[14,3] = "1638"
$data[15,0] = "749, 749, 781, 1074"
$data[15,1] = 1
$data[15,2] = "781"
$data[15,3] = "1811"
$data[16,0] = "456, 749, 781, 1074"
$data[16,1] = 1
$data[16,2] = "781"
$data[16,3] = "1842"
$data[17,0] = "99, 131, 749, 1074, SF"
$data[17,1] = 2
$data[17,2] = "131, 131"
$data[17,3] = "2151, 2065"
$data[18,0] = "99, 131, 781, 1074"
$data[18,1] = 1
$data[18,2] = "781"
$data[18,3] = "2079"
$data[19,0] = "131, 781, 781, 1074"
$data[19,1] = 1
$data[19,2] = "781"
$data[19,3] = "1878"
$data[20,0] = "131, 131, 781, 1074"
$data[20,1] = 1
$data[20,2] = "131"
$data[20,3] = "1940"
$data[21,0] = "131, 749, 781, 781"
$data[21,1] = 1
$data[21,2] = "781"
$data[21,3] = "1988"
$data[22,0] = "99, 749, 781, 1074"
$data[22,1] = 3
$data[22,2] = "781, 781, 781"
$data[22,3] = "2343, 2309, 2430"
$data[23,0] = "99, 456, 781, 781"
$data[23,1] = 1
$data[23,2] = "781"
$data[23,3] = "2113"
$data[24,0] = "131, 456, 749, 781"
$data[24,1] = 1
$data[24,2] = "456"
$data[24,3] = "3680"
$data[25,0] = "99, 749, 781, 1106"
$data[25,1] = 2
$data[25,2] = "781, 781"
$data[25,3] = "2241, 5636"
$data[26,0] = "99, 781, 1074"
$data[26,1] = 1
$data[26,2] = "781"
$data[26,3] = "2833"
$data[27,0] = "131, 131, 781, 1106"
$data[27,1] = 1
$data[27,2] = "781"
$data[27,3] = "2474"
$data[28,0] = "131, 424, 781, 1106"
$data[28,1] = 1
$data[28,2] = "781"
$data[28,3] = "4010"
$data[29,0] = "131, 131, 749, 1074"
$data[29,1] = 1
$data[29,2] = "131"
$data[29,3] = "2605"
$data[30,0] = "99, 749, 781"
$data[30,1] = 1
$data[30,2] = "781"
$data[30,3] = "3008"
$data[31,0] = "131, 749, 781"
$data[31,1] = 2
$data[31,2] = "781, 781"
$data[31,3] = "3280, 3394"
$data[32,0] = "131, 456, 1106"
$data[32,1] = 3
$data[32,2] = "456, 131, 1106"
$data[32,3] = "3861, 4290, 5791"
$data[33,0] = "131, 131, 1074, 1106"
$data[33,1] = 1
$data[33,2] = "131"
$data[33,3] = "3903"
$data[34,0] = "749, 781, 1074"
$data[34,1] = 1
$data[34,2] = "781"
$data[34,3] = "4719"
$data[35,0] = "131, 424, 456, 749"
$data[35,1] = 1
$data[35,2] = "456"
$data[35,3] = "4047"
$data[36,0] = "131, 131, 1106, 1106"
$data[36,1] = 1
$data[36,2] = "131"
$data[36,3] = "4907"
$data[37,0] = "131, 749, 1074, 1106"
$data[37,1] = 1
$data[37,2] = "131"
$data[37,3] = "4062"
$data[38,0] = "456, 749, 1106"
$data[38,1] = 1
$data[38,2] = "456"
$data[38,3] = "4414"
$data[39,0] = "131, 424, 781, 1074"
$data[39,1] = 1
$data[39,2] = "781"
$data[39,3] = "4508"
$data[40,0] = "99, 424, 456, 781, 1074"
$data[40,1] = 1
$data[40,2] = "456"
$data[40,3] = "5470"
$data[41,0] = "99, 456, 749, 1074"
$data[41,1] = 1
$data[41,2] = "456"
$data[41,3] = "4624"
$data[42,0] = "131, 424, 456, 749, 1074"
$data[42,1] = 1
$data[42,2] = "456"
$data[42,3] = "4577"
$data[43,0] = "131, 131, 749, 1106"
$data[43,1] = 1
$data[43,2] = "131"
$data[43,3] = "4695"
$data[44,0] = "99, 99, 749"
$data[44,1] = 1
$data[44,2] = "99"
$data[44,3] = "5704"
$data[45,0] = "99, 456, 749, 781, 1106, SF"
$data[45,1] = 1
$data[45,2] = "456"
$data[45,3] = "5100"
$data[46,0] = "99, 781, 1074, 1074"
$data[46,1] = 1
$data[46,2] = "781"
$data[46,3] = "5275"
$data[47,0] = "99, 424, 456, 456, 749, 1074"
$data[47,1] = 1
$data[47,2] = "456"
$data[47,3] = "4991"
$data[48,0] = "99, 424, 456, 749, 781, 1074"
$data[48,1] = 1
$data[48,2] = "456"
$data[48,3] = "5030"
$data[49,0] = "131, 749, 781, 1074, 1106"
$data[49,1] = 1
$data[49,2] = "1106"
$data[49,3] = "5124"
$data[50,0] = "99, 131, 424, 456, 1074"
$data[50,1] = 1
$data[50,2] = "456"
$data[50,3] = "5223"
$data[51,0] = "99, 424, 781, 1074"
$data[51,1] = 1
$data[51,2] = "781"
$data[51,3] = "5477"
$data[52,0] = "424, 456"
$data[52,1] = 1
$data[52,2] = "456"
$data[52,3] = "5858"
$data[53,0] = "131, 1106"
$data[53,1] = 1
$data[53,2] = "1106"
$data[53,3] = "5837"
$data[54,0] = "131"
$data[54,1] = 1
$data[54,2] = "131"
$data[54,3] = "5872"
$data[55,0] = "781, 1106"
$data[55,1] = 1
$data[55,2] = "781"
$data[55,3] = "6040"
$data[56,0] = "456, 456"
$data[56,1] = 1
$data[56,2] = "456"
$data[56,3] = "5964"
$data[57,0] = "424, 749, 1106"
$data[57,1] = 3
$data[57,2] = "1106, 424, 749"
$data[57,3] = "6798, 6798, 6798"
$data[58,0] = "456, 1074"
$data[58,1] = 1
$data[58,2] = "456"
$data[58,3] = "6080"
$data[59,0] = "456, 456, 1074"
$data[59,1] = 1
$data[59,2] = "456"
$data[59,3] = "6152"
$data[60,0] = "99, 456, 1074"
$data[60,1] = 1
$data[60,2] = "456"
$data[60,3] = "6176"
$data[61,0] = "424, 456, 1074"
$data[61,1] = 1
$data[61,2] = "456"
$data[61,3] = "6231"
$data[62,0] = "424, 456, 749, 1106, SF"
$data[62,1] = 1
$data[62,2] = "456"
$data[62,3] = "6334"
$data[63,0] = "456, 456, 781, 1074"
$data[63,1] = 1
$data[63,2] = "456"
$data[63,3] = "6452"
$data[64,0] = "424, 456, 749, 1074"
$data[64,1] = 1
$data[64,2] = "456"
$data[64,3] = "6466"
$data[65,0] = "424, 749, 781, 1106"
$data[65,1] = 2
$data[65,2] = "781, 781"
$data[65,3] = "6503, 6532"
$data[66,0] = "131, 456, 749, 1074"
$data[66,1] = 1
$data[66,2] = "456"
$data[66,3] = "6582"
$data[67,0] = "99, 131, 456, 1106"
$data[67,1] = 1
$data[67,2] = "131"
$data[67,3] = "6737"

$ws.Range("A2:D69").Value = $data

